# Fruta / hortaliza, semanal
# Insert two new weekly records for "Vega Monumental Concepción" (Ciruela, Black Amber)
# right after row 69, pushing the existing rows 70-85 down to 72-87.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows starting at row 70 (existing rows shift down by 2).
$ws.Rows.Item(70).Resize(2).Insert()

# New row 70: Black Amber, Primera
$ws.Range("A70").Value = 11
$ws.Range("B70").Value = "Vega Monumental Concepción"
$ws.Range("C70").Value = "Bíobío"
$ws.Range("D70").Value = 44946
$ws.Range("E70").Value = 8
$ws.Range("F70").Value = "Fruta"
$ws.Range("G70").Value = 100103
$ws.Range("H70").Value = "Frutos de hueso (carozo)"
$ws.Range("I70").Value = 100103002
$ws.Range("J70").Value = "Ciruela"
$ws.Range("K70").Value = "Black Amber"
$ws.Range("L70").Value = "Primera"
$ws.Range("M70").Value = 100
$ws.Range("N70").Value = 11000
$ws.Range("O70").Value = 12000
$ws.Range("P70").Value = 11500
$ws.Range("Q70").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R70").Value = "Región de O'Higgins"
$ws.Range("S70").Value = 639
$ws.Range("T70").Value = 18

# New row 71: Black Amber, Segunda
$ws.Range("A71").Value = 11
$ws.Range("B71").Value = "Vega Monumental Concepción"
$ws.Range("C71").Value = "Bíobío"
$ws.Range("D71").Value = 44946
$ws.Range("E71").Value = 8
$ws.Range("F71").Value = "Fruta"
$ws.Range("G71").Value = 100103
$ws.Range("H71").Value = "Frutos de hueso (carozo)"
$ws.Range("I71").Value = 100103002
$ws.Range("J71").Value = "Ciruela"
$ws.Range("K71").Value = "Black Amber"
$ws.Range("L71").Value = "Segunda"
$ws.Range("M71").Value = 50
$ws.Range("N71").Value = 9000
$ws.Range("O71").Value = 9000
$ws.Range("P71").Value = 9000
$ws.Range("Q71").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R71").Value = "Región de O'Higgins"
$ws.Range("S71").Value = 500
$ws.Range("T71").Value = 18
